# Generate Report for Handback
# Adds a new handed-back file entry (8e412902-250a-4c9e-a465-d970d91d88b4)
# as row 4 on the "Overview", "zh-cn" and "de-de" worksheets.

$wb = $excel.ActiveWorkbook

$uuid = "8e412902-250a-4c9e-a465-d970d91d88b4"
$hash = "53602b9c5e3ce60059b767b48d1c9b439d283269"
$statusInSync = "Handed back: in sync with en-US"
$hyperlinkColor = 15570276   # matches the workbook's custom HyperLink font color (#FF6495ED)

function Style-AsLink($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "$uuid.md"
$wsOverview.Range("B4").Value = $statusInSync
$wsOverview.Range("C4").Value = $statusInSync

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/30877432d1026706d7e805da846a32c3bb81e3c2/e2e/$uuid.md",
    "",
    "",
    "$uuid.md") | Out-Null
Style-AsLink $wsOverview.Range("A4")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A4").Value = "$uuid.md"
$wsZhCn.Range("B4").Value = $statusInSync
$wsZhCn.Range("C4").Value = "$uuid.$hash.zh-cn.xlf"
$wsZhCn.Range("D4").Value = "2016-02-18 03:30:01"
$wsZhCn.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("E4").Value = "$uuid.md"
$wsZhCn.Range("F4").Value = "$uuid.$hash.zh-cn.xlf"
$wsZhCn.Range("G4").Value = "2016-02-18 03:30:44"
$wsZhCn.Range("H4").Value = "Include"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/30877432d1026706d7e805da846a32c3bb81e3c2/e2e/$uuid.md",
    "",
    "",
    "$uuid.md") | Out-Null
Style-AsLink $wsZhCn.Range("A4")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9b62179273c8eb5bb682575ec87a171ac826a6fc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$uuid.$hash.zh-cn.xlf",
    "",
    "",
    "$uuid.$hash.zh-cn.xlf") | Out-Null
Style-AsLink $wsZhCn.Range("C4")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/e48478dcb74f21345d2cce8038a39d5e0853964b/e2e/$uuid.md",
    "",
    "",
    "$uuid.md") | Out-Null
Style-AsLink $wsZhCn.Range("E4")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/50af03b971722f244f58d669cbee3772a0770217/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$uuid.$hash.zh-cn.xlf",
    "",
    "",
    "$uuid.$hash.zh-cn.xlf") | Out-Null
Style-AsLink $wsZhCn.Range("F4")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A4").Value = "$uuid.md"
$wsDeDe.Range("B4").Value = $statusInSync
$wsDeDe.Range("C4").Value = "$uuid.$hash.de-de.xlf"
$wsDeDe.Range("D4").Value = "2016-02-18 03:30:14"
$wsDeDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("E4").Value = "$uuid.md"
$wsDeDe.Range("F4").Value = "$uuid.$hash.de-de.xlf"
$wsDeDe.Range("G4").Value = "2016-02-18 03:31:05"
$wsDeDe.Range("H4").Value = "Include"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTest/oltest/blob/30877432d1026706d7e805da846a32c3bb81e3c2/e2e/$uuid.md",
    "",
    "",
    "$uuid.md") | Out-Null
Style-AsLink $wsDeDe.Range("A4")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("C4"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/21a278f64f7fd633dbdde131ca3766e4d58e72e3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$uuid.$hash.de-de.xlf",
    "",
    "",
    "$uuid.$hash.de-de.xlf") | Out-Null
Style-AsLink $wsDeDe.Range("C4")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("E4"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/10275dff6c15c0c8e9df469611a11f5125227c37/e2e/$uuid.md",
    "",
    "",
    "$uuid.md") | Out-Null
Style-AsLink $wsDeDe.Range("E4")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F4"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/12da86a78c49ea20e32684b27b95e90934833489/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$uuid.$hash.de-de.xlf",
    "",
    "",
    "$uuid.$hash.de-de.xlf") | Out-Null
Style-AsLink $wsDeDe.Range("F4")

Write-Host "Done adding handback row for $uuid"
